$d = $word.ActiveDocument

# 1) Remove the introductory paragraph entirely ("If you're not doing this
#    module straight after the last one, ... create the data frame 'data'.")
$introPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "If you*re not doing this module*") {
        $introPara = $p
        break
    }
}
if ($introPara -ne $null) {
    $introPara.Range.Delete()
}

# 2) Rework the "Summary is a useful function..." paragraph: split the " line"
#    run and wrap "line" in proofErr gramStart/gramEnd, wrap the "1st" value in
#    gramStart/gramEnd, and split "columns with categorical data it will " so
#    "data" is wrapped in gramStart/gramEnd as well.
$summaryXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Summary is a </w:t></w:r><w:r><w:t>useful function for giving you an overview of your data. If we run that</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>line</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> we can see </w:t></w:r><w:r><w:t xml:space="preserve">the output </w:t></w:r><w:r><w:t>down in the console</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve"> For the columns which are numeric we have </w:t></w:r><w:r><w:t xml:space="preserve">the minimum and maximum values, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>1</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>st</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and 3</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>rd</w:t></w:r><w:r><w:t xml:space="preserve"> quartile, median and mean</w:t></w:r><w:r><w:t xml:space="preserve">. For </w:t></w:r><w:r><w:t xml:space="preserve">columns with categorical </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>data</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> it will </w:t></w:r><w:r><w:t>show the count for the first few rows and then a total for how many are in other categories</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> If there were any </w:t></w:r><w:r><w:t xml:space="preserve">missing values </w:t></w:r><w:r><w:t xml:space="preserve">in our dataset </w:t></w:r><w:r><w:t>it would tell us that as well.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$summaryPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Summary is a useful function*") {
        $summaryPara = $p
        break
    }
}
if ($summaryPara -ne $null) {
    $rng = $d.Range($summaryPara.Range.Start, $summaryPara.Range.End)
    $rng.InsertXML($summaryXml)
}

# 3) Merge the trailing ")" run with the " with `` around any name with spaces
#    included." run into a single run.
$answerXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Answer = c</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">The function we need is </w:t></w:r><w:r><w:t>rename(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>new_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>old_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>) with `` around any name with spaces included.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$answerPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Answer = c*") {
        $answerPara = $p
        break
    }
}
if ($answerPara -ne $null) {
    $rng = $d.Range($answerPara.Range.Start, $answerPara.Range.End)
    $rng.InsertXML($answerXml)
}

Write-Output "DONE"
foreach ($p in $d.Paragraphs) {
    Write-Output $p.Range.Text
}
